$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Baza podataka" -> "Analiza sadržaja" in cell B2
$ws.Range("B2").Value = "Analiza sadržaja"

# Update the selection to reflect the edited cell
$ws.Range("B2").Select()
